# Quarterly indexing esoteric bug-fix operation
#
# Column A held the first day of each quarter (e.g. 1988-07-01) as a date
# serial. The corrected series instead stamps each quarterly observation
# on the 15th of the *second* month of the quarter (e.g. 1988-08-15) -
# i.e. one month plus fourteen days after the old quarter-start date.
# Recompute every date in column A (rows 2-150) from its own prior value
# so the fix is self-contained and not dependent on hard-coded constants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 150) { $lastRow = 150 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    if ($serial -eq $null) { continue }

    $oldDate = [DateTime]::FromOADate($serial)

    $newYear = $oldDate.Year
    $newMonth = $oldDate.Month + 1
    if ($newMonth -gt 12) {
        $newMonth = $newMonth - 12
        $newYear = $newYear + 1
    }

    $newDate = Get-Date -Year $newYear -Month $newMonth -Day 15 -Hour 0 -Minute 0 -Second 0

    $cell.Value = $newDate
}
